$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.596.54'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '2.507.00'
$ws.Range('E3').Value = '  -1.75%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'311.57"
$ws.Range('E5').Value = '  +3.04%  '
$ws.Range('D6').Value = "'95.82"
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('D7').Value = "'0.585"
$ws.Range('E7').Value = '  +1.98%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.541"
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('D10').Value = "'36.03"
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').Value = "'0.0813"
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = "'7.74"
$ws.Range('E12').Value = '  +2.97%  '
$ws.Range('D13').Value = "'0.113"
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D14').Value = '2.892.16'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').Value = "'15.60"
$ws.Range('E15').Value = '  +7.86%  '
$ws.Range('D16').Value = '2.497.86'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = "'0.857"
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '42.592.61'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').Value = "'12.89"
$ws.Range('E19').Value = '  -3.91%  '
$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').Value = "'71.52"
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = "'253.61"
$ws.Range('E23').Value = '  -1.02%  '
$ws.Range('D24').Value = "'2.97"
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('E25').Value = '  -2.36%  '
$ws.Range('D26').Value = "'27.07"
$ws.Range('E26').Value = '  -3.57%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = "'2.35"
$ws.Range('E28').Value = '  +11.88%  '
$ws.Range('D29').Value = "'10.15"
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').Value = "'37.75"
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('D31').Value = "'5.93"
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = "'154.19"
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').Value = "'19.29"
$ws.Range('E33').Value = '  +5.26%  '
$ws.Range('D34').Value = "'3.30"
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'2.08"
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.0786"
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('E37').Value = '  -4.66%  '
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').Value = "'24.59"
$ws.Range('E39').Value = '  -4.74%  '
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = "'3.89"
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = "'2.03"
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').Value = "'0.0303"
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'0.999"
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '2.024.34'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').Value = "'84.37"
$ws.Range('E47').Value = '  -3.71%  '
$ws.Range('D48').Value = "'8.95"
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').Value = '2.750.88'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = "'73.17"
$ws.Range('E50').Value = '  -3.48%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.191"
$ws.Range('E51').Value = '  +0.89%  '
